$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NATMI LR-pair metrics with new TPM-based expression values.
# Source (raw) values changed per the new TPM recompute:
#   Ligand avg expr (G) by sending cluster: ECs 0.05916633333333333 -> 1.814087666666667;
#                                           MuSCs 4.322670666666666 -> 7.416845666666667 (FAPs unchanged)
#   Ligand-expressing cells (E) ECs: 1 -> 3 (FAPs/MuSCs unchanged)
#   Receptor avg expr (M) by target cluster: ECs 0.502378 -> 1.837384;
#     Inflammatory-Mac 5.474800666666667 -> 5.528959666666668; MuSCs 7.461044333333334 -> 3.046454666666667;
#     Neutrophils 7.423250333333333 -> 9.021246333333332; Resolving-Mac 5.012936333333333 -> 8.883473666666667 (FAPs unchanged)
# All dependent columns (detection rate, totals, specificities, edge weights) are updated to match.

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.814087666666667
$ws.Range("H2").Value = 5.442263
$ws.Range("I2").Value = 0.1211063206477811
$ws.Range("J2").Value = 0.1211063206477811
$ws.Range("M2").Value = 1.837384
$ws.Range("N2").Value = 5.512152
$ws.Range("O2").Value = 0.0635335947613339
$ws.Range("P2").Value = 0.0635335947613339
$ws.Range("Q2").Value = 3.333175653330667
$ws.Range("R2").Value = 29.998580879976
$ws.Range("S2").Value = 0.00769431989907229
$ws.Range("T2").Value = 0.00769431989907229

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.814087666666667
$ws.Range("H3").Value = 5.442263
$ws.Range("I3").Value = 0.1211063206477811
$ws.Range("J3").Value = 0.1211063206477811
$ws.Range("O3").Value = 0.02082867030699976
$ws.Range("P3").Value = 0.02082867030699976
$ws.Range("Q3").Value = 1.092738684460444
$ws.Range("R3").Value = 9.834648160143999
$ws.Range("S3").Value = 0.002522483624866431
$ws.Range("T3").Value = 0.002522483624866431

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.814087666666667
$ws.Range("H4").Value = 5.442263
$ws.Range("I4").Value = 0.1211063206477811
$ws.Range("J4").Value = 0.1211063206477811
$ws.Range("M4").Value = 5.528959666666668
$ws.Range("N4").Value = 16.586879
$ws.Range("O4").Value = 0.1911819646376369
$ws.Range("P4").Value = 0.1911819646376369
$ws.Range("Q4").Value = 10.03001754079745
$ws.Range("R4").Value = 90.27015786717701
$ws.Range("S4").Value = 0.0231533443114784
$ws.Range("T4").Value = 0.0231533443114784

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.814087666666667
$ws.Range("H5").Value = 5.442263
$ws.Range("I5").Value = 0.1211063206477811
$ws.Range("J5").Value = 0.1211063206477811
$ws.Range("M5").Value = 3.046454666666667
$ws.Range("N5").Value = 9.139364
$ws.Range("O5").Value = 0.1053411895666744
$ws.Range("P5").Value = 0.1053411895666744
$ws.Range("Q5").Value = 5.526535837859111
$ws.Range("R5").Value = 49.738822540732
$ws.Range("S5").Value = 0.01275748388108037
$ws.Range("T5").Value = 0.01275748388108037

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.814087666666667
$ws.Range("H6").Value = 5.442263
$ws.Range("I6").Value = 0.1211063206477811
$ws.Range("J6").Value = 0.1211063206477811
$ws.Range("M6").Value = 9.021246333333332
$ws.Range("N6").Value = 27.063739
$ws.Range("O6").Value = 0.3119392618985303
$ws.Range("P6").Value = 0.3119392618985303
$ws.Range("Q6").Value = 16.36533171126188
$ws.Range("R6").Value = 147.287985401357
$ws.Range("S6").Value = 0.03777781627411559
$ws.Range("T6").Value = 0.03777781627411559

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.814087666666667
$ws.Range("H7").Value = 5.442263
$ws.Range("I7").Value = 0.1211063206477811
$ws.Range("J7").Value = 0.1211063206477811
$ws.Range("M7").Value = 8.883473666666667
$ws.Range("N7").Value = 26.650421
$ws.Range("O7").Value = 0.3071753188288246
$ws.Range("P7").Value = 0.3071753188288246
$ws.Range("Q7").Value = 16.11540001585811
$ws.Range("R7").Value = 145.038600142723
$ws.Range("S7").Value = 0.03720087265716803
$ws.Range("T7").Value = 0.03720087265716803

# Row 8
$ws.Range("I8").Value = 0.3837539427192561
$ws.Range("J8").Value = 0.3837539427192561
$ws.Range("M8").Value = 1.837384
$ws.Range("N8").Value = 5.512152
$ws.Range("O8").Value = 0.0635335947613339
$ws.Range("P8").Value = 0.0635335947613339
$ws.Range("Q8").Value = 10.56195326469867
$ws.Range("R8").Value = 95.05757938228801
$ws.Range("S8").Value = 0.02438126748478936
$ws.Range("T8").Value = 0.02438126748478936

# Row 9
$ws.Range("I9").Value = 0.3837539427192561
$ws.Range("J9").Value = 0.3837539427192561
$ws.Range("O9").Value = 0.02082867030699976
$ws.Range("P9").Value = 0.02082867030699976
$ws.Range("S9").Value = 0.007993084351910654
$ws.Range("T9").Value = 0.007993084351910654

# Row 10
$ws.Range("I10").Value = 0.3837539427192561
$ws.Range("J10").Value = 0.3837539427192561
$ws.Range("M10").Value = 5.528959666666668
$ws.Range("N10").Value = 16.586879
$ws.Range("O10").Value = 0.1911819646376369
$ws.Range("P10").Value = 0.1911819646376369
$ws.Range("Q10").Value = 31.78247639129179
$ws.Range("R10").Value = 286.0422875216261
$ws.Range("S10").Value = 0.07336683270650654
$ws.Range("T10").Value = 0.07336683270650654

# Row 11
$ws.Range("I11").Value = 0.3837539427192561
$ws.Range("J11").Value = 0.3837539427192561
$ws.Range("M11").Value = 3.046454666666667
$ws.Range("N11").Value = 9.139364
$ws.Range("O11").Value = 0.1053411895666744
$ws.Range("P11").Value = 0.1053411895666744
$ws.Range("Q11").Value = 17.51213236446845
$ws.Range("R11").Value = 157.609191280216
$ws.Range("S11").Value = 0.04042509682694788
$ws.Range("T11").Value = 0.04042509682694788

# Row 12
$ws.Range("I12").Value = 0.3837539427192561
$ws.Range("J12").Value = 0.3837539427192561
$ws.Range("M12").Value = 9.021246333333332
$ws.Range("N12").Value = 27.063739
$ws.Range("O12").Value = 0.3119392618985303
$ws.Range("P12").Value = 0.3119392618985303
$ws.Range("Q12").Value = 51.85741367182955
$ws.Range("R12").Value = 466.716723046466
$ws.Range("S12").Value = 0.1197079216424956
$ws.Range("T12").Value = 0.1197079216424956

# Row 13
$ws.Range("I13").Value = 0.3837539427192561
$ws.Range("J13").Value = 0.3837539427192561
$ws.Range("M13").Value = 8.883473666666667
$ws.Range("N13").Value = 26.650421
$ws.Range("O13").Value = 0.3071753188288246
$ws.Range("P13").Value = 0.3071753188288246
$ws.Range("Q13").Value = 51.06544614273045
$ws.Range("R13").Value = 459.5890152845741
$ws.Range("S13").Value = 0.117879739706606
$ws.Range("T13").Value = 0.117879739706606

# Row 14
$ws.Range("G14").Value = 7.416845666666667
$ws.Range("H14").Value = 22.250537
$ws.Range("I14").Value = 0.4951397366329628
$ws.Range("J14").Value = 0.4951397366329628
$ws.Range("M14").Value = 1.837384
$ws.Range("N14").Value = 5.512152
$ws.Range("O14").Value = 0.0635335947613339
$ws.Range("P14").Value = 0.0635335947613339
$ws.Range("Q14").Value = 13.62759355840267
$ws.Range("R14").Value = 122.648342025624
$ws.Range("S14").Value = 0.03145800737747225
$ws.Range("T14").Value = 0.03145800737747225

# Row 15
$ws.Range("G15").Value = 7.416845666666667
$ws.Range("H15").Value = 22.250537
$ws.Range("I15").Value = 0.4951397366329628
$ws.Range("J15").Value = 0.4951397366329628
$ws.Range("O15").Value = 0.02082867030699976
$ws.Range("P15").Value = 0.02082867030699976
$ws.Range("Q15").Value = 4.467630934028445
$ws.Range("R15").Value = 40.208678406256
$ws.Range("S15").Value = 0.01031310233022267
$ws.Range("T15").Value = 0.01031310233022267

# Row 16
$ws.Range("G16").Value = 7.416845666666667
$ws.Range("H16").Value = 22.250537
$ws.Range("I16").Value = 0.4951397366329628
$ws.Range("J16").Value = 0.4951397366329628
$ws.Range("M16").Value = 5.528959666666668
$ws.Range("N16").Value = 16.586879
$ws.Range("O16").Value = 0.1911819646376369
$ws.Range("P16").Value = 0.1911819646376369
$ws.Range("Q16").Value = 41.00744054489146
$ws.Range("R16").Value = 369.0669649040231
$ws.Range("S16").Value = 0.09466178761965195
$ws.Range("T16").Value = 0.09466178761965195

# Row 17
$ws.Range("G17").Value = 7.416845666666667
$ws.Range("H17").Value = 22.250537
$ws.Range("I17").Value = 0.4951397366329628
$ws.Range("J17").Value = 0.4951397366329628
$ws.Range("M17").Value = 3.046454666666667
$ws.Range("N17").Value = 9.139364
$ws.Range("O17").Value = 0.1053411895666744
$ws.Range("P17").Value = 0.1053411895666744
$ws.Range("Q17").Value = 22.59508409316312
$ws.Range("R17").Value = 203.355756838468
$ws.Range("S17").Value = 0.0521586088586462
$ws.Range("T17").Value = 0.0521586088586462

# Row 18
$ws.Range("G18").Value = 7.416845666666667
$ws.Range("H18").Value = 22.250537
$ws.Range("I18").Value = 0.4951397366329628
$ws.Range("J18").Value = 0.4951397366329628
$ws.Range("M18").Value = 9.021246333333332
$ws.Range("N18").Value = 27.063739
$ws.Range("O18").Value = 0.3119392618985303
$ws.Range("P18").Value = 0.3119392618985303
$ws.Range("Q18").Value = 66.90919177531589
$ws.Range("R18").Value = 602.182725977843
$ws.Range("S18").Value = 0.1544535239819191
$ws.Range("T18").Value = 0.1544535239819191

# Row 19
$ws.Range("G19").Value = 7.416845666666667
$ws.Range("H19").Value = 22.250537
$ws.Range("I19").Value = 0.4951397366329628
$ws.Range("J19").Value = 0.4951397366329628
$ws.Range("M19").Value = 8.883473666666667
$ws.Range("N19").Value = 26.650421
$ws.Range("O19").Value = 0.3071753188288246
$ws.Range("P19").Value = 0.3071753188288246
$ws.Range("Q19").Value = 65.88735316956412
$ws.Range("R19").Value = 592.9861785260771
$ws.Range("S19").Value = 0.1520947064650506
$ws.Range("T19").Value = 0.1520947064650506
